$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "position" column (I) header
$ws.Range("I1").Value = "position"

# Rows 2-10  -> pns
for ($r = 2; $r -le 10; $r++) {
    $ws.Cells.Item($r, 9).Value = "pns"
}

# Rows 11-20 -> p3k
for ($r = 11; $r -le 20; $r++) {
    $ws.Cells.Item($r, 9).Value = "p3k"
}

# Rows 21-29 -> cpns
for ($r = 21; $r -le 29; $r++) {
    $ws.Cells.Item($r, 9).Value = "cpns"
}

# Match text-formatted style (s="1") used by the rest of the table
$ws.Range("I1:I29").NumberFormat = "@"

# Size the new column (closest achievable match to the source's 19.140625 char width)
$ws.Columns.Item(9).ColumnWidth = 18.1

# Move / extend the active selection like the source workbook
$null = $ws.Range("I21:I29").Select()
